$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 data refresh (PERIOD TO EXPIRE / LAST UPDATE) ---------------
# H3: PERIOD TO EXPIRE -47 -> -55
$ws.Range("H3").Value = -55

# I3: LAST UPDATE 08-Sep-2025 -> 16-Sep-2025 (kept as literal text, not
# an auto-converted date serial, by forcing the cell to Text format
# before typing the new value).
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"

# --- Header styling: title + column headers get a bold white font ------
# (previously the title used a bold 14pt font and the column headers a
# plain bold font with no color; both now share one bold, default-size,
# white font)
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215

$ws.Range("A2:K2").Font.Bold = $true
$ws.Range("A2:K2").Font.Color = 16777215
